# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 4-7 hold F (想去人数) values that changed
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1826
$wsExhibit.Range("F5").Value = 805
$wsExhibit.Range("F6").Value = 509
$wsExhibit.Range("F7").Value = 220

# Sheet "全部类型" - rows 4, 6, 7, 8 hold the same items' F (想去人数) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1826
$wsAll.Range("F6").Value = 805
$wsAll.Range("F7").Value = 509
$wsAll.Range("F8").Value = 220
